$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously held Graph-search timings (cols A:F, rows 2:6).
# Replace that with the new Tree-search parallel timings (cols C:D and G:J,
# rows 2:4) per "display parallel data for tree".
$ws.Range("A2:L6").ClearContents()

# Row 2 - Graph BFS/DFS 1000 cols (C2:D2) + Tree BFS/DFS 10000 (G2:H2) + Tree BFS/DFS 1000 (I2:J2)
$ws.Range("C2").Value = 0.0207863
$ws.Range("D2").Value = 0.0256373
$ws.Range("G2").Value = 0.0001315
$ws.Range("H2").Value = 0.0003602
$ws.Range("I2").Value = 0.0004385
$ws.Range("J2").Value = 0.0001274

# Row 3 - Tree BFS/DFS 10000 (G3:H3) + Tree BFS/DFS 1000 (I3:J3)
$ws.Range("G3").Value = 0.0008679
$ws.Range("H3").Value = 0.0010293
$ws.Range("I3").Value = 0.0003156
$ws.Range("J3").Value = 0.0002806

# Row 4 - Tree BFS/DFS 1000 (I4:J4)
$ws.Range("I4").Value = 0.0003374
$ws.Range("J4").Value = 0.0005466

# Update the selected cell shown when the workbook is reopened
$ws.Range("E11").Select()
